# Fix contact information missing from short resumes:
# insert a centered "contact info" paragraph right after the name
# heading ("Dheeraj Chand") and before "PROFESSIONAL SUMMARY".

$d = $word.ActiveDocument

# The name line is the document's first paragraph; restrict the
# Find/Replace to its Range so we only touch that exact text.
$nameRange = $d.Paragraphs(1).Range

$contact = "202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX"

# Replace "Dheeraj Chand" with "Dheeraj Chand" + a new paragraph mark
# (^p) + the contact-info text. The new paragraph inherits the
# surrounding run/paragraph context from the Find/Replace operation
# rather than copying the bold/28pt name styling, matching how the
# contact line should render: plain text, centered.
$nameRange.Find.Execute(
    "Dheeraj Chand",  # FindText
    $true,            # MatchCase
    $true,            # MatchWholeWord
    $false,           # MatchWildcards
    $false,           # MatchSoundsLike
    $false,           # MatchAllWordForms
    $true,            # Forward
    1,                # Wrap (wdFindContinue)
    $false,           # Format
    "Dheeraj Chand^p$contact",  # ReplaceWith
    2                 # Replace (wdReplaceAll)
)
